$d = $word.ActiveDocument

$pairs = @(
    @("53÷6=", "49÷2="),
    @("85÷5=", "58÷4="),
    @("49÷5=", "62÷3="),
    @("35÷8=", "77÷3="),
    @("84÷6=", "28÷6="),
    @("21÷2=", "14÷3="),
    @("69÷5=", "63÷6="),
    @("52÷6=", "23÷2="),
    @("80÷7=", "54÷4="),
    @("32÷9=", "64÷9="),
    @("73÷4=", "15÷5="),
    @("89÷7=", "15÷4="),
    @("54÷3=", "28÷4="),
    @("38÷5=", "28÷2="),
    @("99÷9=", "19÷5="),
    @("85÷9=", "30÷2="),
    @("62÷8=", "83÷9="),
    @("93÷5=", "51÷5="),
    @("75÷2=", "76÷2="),
    @("48÷5=", "29÷4="),
    @("14÷5=", "50÷2="),
    @("61÷3=", "69÷3="),
    @("80÷4=", "38÷7="),
    @("11÷8=", "39÷4="),
    @("58÷7=", "40÷2=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
